$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'310.53"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'-0.42%"
$ws.Range("E2").ClearFormats()
$ws.Range("G2").Value = "'6"
$ws.Range("G2").ClearFormats()

$ws.Range("D3").Value = "'37.69"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'0.39%"
$ws.Range("E3").ClearFormats()
$ws.Range("G3").Value = "'6"
$ws.Range("G3").ClearFormats()

$ws.Range("D4").Value = "'5.170"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'2.09%"
$ws.Range("E4").ClearFormats()
$ws.Range("G4").Value = "'6"
$ws.Range("G4").ClearFormats()

$ws.Range("D5").Value = "'0.07904"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'1.98%"
$ws.Range("E5").ClearFormats()
$ws.Range("G5").Value = "'6"
$ws.Range("G5").ClearFormats()

$ws.Range("D6").Value = "'1.910"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'1.51%"
$ws.Range("E6").ClearFormats()
$ws.Range("G6").Value = "'6"
$ws.Range("G6").ClearFormats()

$ws.Range("D7").Value = "'8.268"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'0.94%"
$ws.Range("E7").ClearFormats()
$ws.Range("G7").Value = "'6"
$ws.Range("G7").ClearFormats()

$ws.Range("D8").Value = "'2.999"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'1.18%"
$ws.Range("E8").ClearFormats()
$ws.Range("G8").Value = "'6"
$ws.Range("G8").ClearFormats()

$ws.Range("D9").Value = "'0.9382"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'2.25%"
$ws.Range("E9").ClearFormats()
$ws.Range("G9").Value = "'6"
$ws.Range("G9").ClearFormats()

$ws.Range("D10").Value = "'0.1112"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'-9.19%"
$ws.Range("E10").ClearFormats()
$ws.Range("G10").Value = "'6"
$ws.Range("G10").ClearFormats()

$ws.Range("D11").Value = "'0.1949"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'2.43%"
$ws.Range("E11").ClearFormats()
$ws.Range("G11").Value = "'6"
$ws.Range("G11").ClearFormats()

$ws.Range("D12").Value = "'0.09099"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'2.46%"
$ws.Range("E12").ClearFormats()
$ws.Range("G12").Value = "'6"
$ws.Range("G12").ClearFormats()

$ws.Range("D13").Value = "'0.03332"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'-1.65%"
$ws.Range("E13").ClearFormats()
$ws.Range("G13").Value = "'6"
$ws.Range("G13").ClearFormats()

$ws.Range("D14").Value = "'0.09609"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'-0.94%"
$ws.Range("E14").ClearFormats()
$ws.Range("G14").Value = "'6"
$ws.Range("G14").ClearFormats()

$ws.Range("D15").Value = "'0.001393"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'1.01%"
$ws.Range("E15").ClearFormats()
$ws.Range("G15").Value = "'6"
$ws.Range("G15").ClearFormats()

$ws.Range("D16").Value = "'0.005729"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'-2.31%"
$ws.Range("E16").ClearFormats()
$ws.Range("G16").Value = "'6"
$ws.Range("G16").ClearFormats()

$ws.Range("D17").Value = "'3.595"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'1.55%"
$ws.Range("E17").ClearFormats()
$ws.Range("G17").Value = "'6"
$ws.Range("G17").ClearFormats()

$ws.Range("D18").Value = "'4.429"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'1.47%"
$ws.Range("E18").ClearFormats()
$ws.Range("G18").Value = "'6"
$ws.Range("G18").ClearFormats()

$ws.Range("D19").Value = "'0.3412"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'0.19%"
$ws.Range("E19").ClearFormats()
$ws.Range("G19").Value = "'6"
$ws.Range("G19").ClearFormats()

$ws.Range("D20").Value = "'6.423"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'27.28%"
$ws.Range("E20").ClearFormats()
$ws.Range("G20").Value = "'6"
$ws.Range("G20").ClearFormats()

$ws.Range("D21").Value = "'0.1282"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'0.03%"
$ws.Range("E21").ClearFormats()
$ws.Range("G21").Value = "'6"
$ws.Range("G21").ClearFormats()

$ws.Range("D22").Value = "'0.2521"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'-2.88%"
$ws.Range("E22").ClearFormats()
$ws.Range("G22").Value = "'6"
$ws.Range("G22").ClearFormats()

$ws.Range("D23").Value = "'0.04399"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'0.20%"
$ws.Range("E23").ClearFormats()
$ws.Range("G23").Value = "'6"
$ws.Range("G23").ClearFormats()

$ws.Range("E24").Value = "'1.42%"
$ws.Range("E24").ClearFormats()
$ws.Range("G24").Value = "'6"
$ws.Range("G24").ClearFormats()

$ws.Range("D25").Value = "'0.004586"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'8.26%"
$ws.Range("E25").ClearFormats()
$ws.Range("G25").Value = "'6"
$ws.Range("G25").ClearFormats()

$ws.Range("E26").Value = "'0.58%"
$ws.Range("E26").ClearFormats()
$ws.Range("G26").Value = "'6"
$ws.Range("G26").ClearFormats()

$ws.Range("D27").Value = "'0.0003992"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'-98.11%"
$ws.Range("E27").ClearFormats()
$ws.Range("G27").Value = "'6"
$ws.Range("G27").ClearFormats()

$ws.Range("G28").Value = "'6"
$ws.Range("G28").ClearFormats()

$ws.Range("G29").Value = "'6"
$ws.Range("G29").ClearFormats()

$ws.Range("G30").Value = "'6"
$ws.Range("G30").ClearFormats()

$ws.Range("G31").Value = "'6"
$ws.Range("G31").ClearFormats()

$ws.Range("G32").Value = "'6"
$ws.Range("G32").ClearFormats()

$ws.Range("G33").Value = "'6"
$ws.Range("G33").ClearFormats()

$ws.Range("G34").Value = "'6"
$ws.Range("G34").ClearFormats()

$ws.Range("G35").Value = "'6"
$ws.Range("G35").ClearFormats()

$ws.Range("G36").Value = "'6"
$ws.Range("G36").ClearFormats()

$ws.Range("G37").Value = "'6"
$ws.Range("G37").ClearFormats()

$ws.Range("G38").Value = "'6"
$ws.Range("G38").ClearFormats()

$ws.Range("D39").Value = "'0.02236"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'6.59%"
$ws.Range("E39").ClearFormats()
$ws.Range("G39").Value = "'6"
$ws.Range("G39").ClearFormats()

$ws.Range("D40").Value = "'0.05128"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'3.93%"
$ws.Range("E40").ClearFormats()
$ws.Range("G40").Value = "'6"
$ws.Range("G40").ClearFormats()

$ws.Range("D41").Value = "'0.007462"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'-4.51%"
$ws.Range("E41").ClearFormats()
$ws.Range("G41").Value = "'6"
$ws.Range("G41").ClearFormats()

$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").Value = "'0.1356"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'1.12%"
$ws.Range("E42").ClearFormats()
$ws.Range("G42").Value = "'6"
$ws.Range("G42").ClearFormats()

$ws.Range("B43").Value = 'Dexo'
$ws.Range("C43").Value = 'https://coinranking.com/coin/QkL_pl546+dexo-dexo'
$ws.Range("D43").Value = "'0.008754"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'-12.56%"
$ws.Range("E43").ClearFormats()
$ws.Range("G43").Value = "'6"
$ws.Range("G43").ClearFormats()

$ws.Range("D44").Value = "'0.002131"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'6.74%"
$ws.Range("E44").ClearFormats()
$ws.Range("G44").Value = "'6"
$ws.Range("G44").ClearFormats()

$ws.Range("D45").Value = "'0.009324"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'-3.33%"
$ws.Range("E45").ClearFormats()
$ws.Range("G45").Value = "'6"
$ws.Range("G45").ClearFormats()

$ws.Range("D46").Value = "'0.00006622"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'0.59%"
$ws.Range("E46").ClearFormats()
$ws.Range("G46").Value = "'6"
$ws.Range("G46").ClearFormats()

$ws.Range("E47").Value = "'-0.11%"
$ws.Range("E47").ClearFormats()
$ws.Range("G47").Value = "'6"
$ws.Range("G47").ClearFormats()

$ws.Range("E48").Value = "'-6.21%"
$ws.Range("E48").ClearFormats()
$ws.Range("G48").Value = "'6"
$ws.Range("G48").ClearFormats()

$ws.Range("E49").Value = "'-40.72%"
$ws.Range("E49").ClearFormats()
$ws.Range("G49").Value = "'6"
$ws.Range("G49").ClearFormats()

$ws.Range("E50").Value = "'-0.11%"
$ws.Range("E50").ClearFormats()
$ws.Range("G50").Value = "'6"
$ws.Range("G50").ClearFormats()

$ws.Range("E51").Value = "'-0.11%"
$ws.Range("E51").ClearFormats()
$ws.Range("G51").Value = "'6"
$ws.Range("G51").ClearFormats()
